$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 use the same (bold/bordered/centered) style as the
# existing header row, so copy the formatting from H1 before setting values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new "I0" (col I) and "IF" (col J) columns, rows 2-11.
$valuesI = @(10, 6, 8, 5, 5, 6, 7, 6, 5, 9)
$valuesJ = @(10, 7, 9, 6, 6, 6, 8, 6, 5, 9)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $valuesI[$i]
    $ws.Cells.Item($row, 10).Value = $valuesJ[$i]
}
